# Insert a new weekly price record for Achicoria at row 197, pushing the
# existing rows 197-251 down to 198-252 (new dimension A1:R252).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 197..251 down by one, duplicating row 197's formatting for the
# freshly inserted row (matches the s="2" date style on column D).
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new record's data.
$ws.Range("A197").Value = 3
$ws.Range("B197").Value = "Femacal de La Calera"
$ws.Range("C197").Value = "Coquimbo"
$ws.Range("D197").Value = 45211
$ws.Range("E197").Value = 5
$ws.Range("F197").Value = 100112010
$ws.Range("G197").Value = "Achicoria"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 50
$ws.Range("K197").Value = 7000
$ws.Range("L197").Value = 7000
$ws.Range("M197").Value = 7000
$ws.Range("N197").Value = "`$/caja 16 unidades"
$ws.Range("O197").Value = "Provincia de Quillota"
$ws.Range("P197").Value = 438
$ws.Range("Q197").Value = 16
$ws.Range("R197").Value = "Hortaliza"
